$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("eco Pipes")

$ws3.Range("C2").Value = "(1112.7556223528027, 2005.6658755842977, 12.456795843159185)"
$ws3.Range("G2").Value = 582.1549878194318
$ws3.Range("H2").Value = 17.45679584315918
$ws3.Range("M2").Value = 4.286542947838694
$ws3.Range("N2").Value = 22.01425449687209
$ws3.Range("O2").Value = 75.7828351372391
$ws3.Range("S2").Value = 1542710.717721494
$ws3.Range("B3").Value = "(1112.7556223528027, 2005.6658755842977, 12.456795843159185)"
$ws3.Range("C3").Value = "(375.6218616416263, 3051.911036071446, 3.0)"
$ws3.Range("G3").Value = 1279.876770634814
$ws3.Range("H3").Value = -9.456795843159185
$ws3.Range("M3").Value = 5.797278181554547
$ws3.Range("N3").Value = -3.498860808711791
$ws3.Range("O3").Value = 79.2816959459509
$ws3.Range("Q3").Value = "Pipe 3, Pipe 6"
$ws3.Range("B4").Value = "(375.6218616416263, 3051.911036071446, 3.0)"
$ws3.Range("C4").Value = "(144.03445550263132, 2923.579361220609, 5.0)"
$ws3.Range("E4").Value = 450
$ws3.Range("F4").Value = 396.6
$ws3.Range("G4").Value = 264.7748958119118
$ws3.Range("I4").Value = 1100
$ws3.Range("L4").Value = 2.473403000433433
$ws3.Range("M4").Value = 3.467351563656115
$ws3.Range("N4").Value = 5.779289275391436
$ws3.Range("O4").Value = 73.50240667055947
$ws3.Range("Q4").Value = "Pipe 4"
$ws3.Range("R4").Value = 1010
$ws3.Range("S4").Value = 267422.6447700309
$ws3.Range("B5").Value = "(144.03445550263132, 2923.579361220609, 5.0)"
$ws3.Range("C5").Value = "(-945.1143240436759, 3200.8674386767834, 10.0)"
$ws3.Range("G5").Value = 1123.903350776459
$ws3.Range("M5").Value = 16.97998095592162
$ws3.Range("N5").Value = 22.28116865809062
$ws3.Range("O5").Value = 51.22123801246885
$ws3.Range("S5").Value = 1074451.603342295
$ws3.Range("B6").Value = "(-945.1143240436759, 3200.8674386767834, 10.0)"
$ws3.Range("C6").Value = "(-1100.785987141615, 3568.9949130421483, 0.0)"
$ws3.Range("G6").Value = 399.8143370044408
$ws3.Range("M6").Value = 4.214339372476978
$ws3.Range("N6").Value = -5.645555581894074
$ws3.Range("O6").Value = 56.86679359436292
$ws3.Range("S6").Value = 207903.4552423092
$ws3.Range("B7").Value = "(375.6218616416263, 3051.911036071446, 3.0)"
$ws3.Range("C7").Value = "(1370.9059072682867, 3714.0837835015136, 10.0)"
$ws3.Range("D7").Value = "Steel"
$ws3.Range("E7").Value = 12
$ws3.Range("F7").Value = 300
$ws3.Range("G7").Value = 1195.454758206289
$ws3.Range("H7").Value = 7
$ws3.Range("J7").Value = 200
$ws3.Range("L7").Value = 2.35785100876882
$ws3.Range("M7").Value = 19.8393877994091
$ws3.Range("N7").Value = 27.12286023361021
$ws3.Range("O7").Value = 29.74393336075271
$ws3.Range("P7").Value = "Pipe 2"
$ws3.Range("R7").Value = 625
$ws3.Range("S7").Value = 747159.2238789306
$ws3.Range("B8").Value = "(1370.9059072682867, 3714.0837835015136, 10.0)"
$ws3.Range("C8").Value = "(705.2600463000927, 4199.2185612304365, 10.0)"
$ws3.Range("D8").Value = "PE100-16"
$ws3.Range("E8").Value = 315
$ws3.Range("F8").Value = 257.8
$ws3.Range("G8").Value = 823.6747930986961
$ws3.Range("H8").Value = 0
$ws3.Range("I8").Value = 400
$ws3.Range("J8").Value = 100
$ws3.Range("K8").Value = 0
$ws3.Range("L8").Value = 2.128637864434592
$ws3.Range("M8").Value = 13.49778821044812
$ws3.Range("N8").Value = 13.72882528764089
$ws3.Range("O8").Value = 16.01510807311181
$ws3.Range("R8").Value = 690
$ws3.Range("S8").Value = 568335.6072381004
$ws3.Range("B9").Value = "(705.2600463000927, 4199.2185612304365, 10.0)"
$ws3.Range("C9").Value = "(791.7149329602871, 4749.743197126708, 20.0)"
$ws3.Range("D9").Value = "Steel"
$ws3.Range("E9").Value = 10
$ws3.Range("F9").Value = 253
$ws3.Range("G9").Value = 557.36148248345
$ws3.Range("H9").Value = 10
$ws3.Range("I9").Value = 300
$ws3.Range("J9").Value = 300
$ws3.Range("L9").Value = 1.657630886197205
$ws3.Range("M9").Value = 5.875003027482857
$ws3.Range("N9").Value = 16.0151080731118
$ws3.Range("O9").Value = [double]"1.06581410364015E-14"
$ws3.Range("R9").Value = 520
$ws3.Range("S9").Value = 289827.970891394

$ws4 = $wb.Worksheets.Item("eco Pumps")
$ws4.Range("D2").Value = 97.79708963411119
$ws4.Range("G2").Value = 781.2403674475113
